# CDM_Building.xlsx — rework the Sheet1 header row:
#   * reorder / relabel the building-twin attribute columns A1:V1
#     (adds EnergyType / WaterType / HeatingType1 / HeatingType2 /
#     AirConditioning / CertificateId / Guid, drops AreaMeasurementId,
#     SupplyEnergy/SupplyWater/SupplyHeating, GUID)
#   * bold the new header row
#   * drop the stray formatted-but-empty cells at A3 / A5 / A20
#   * widen the columns whose header text no longer fits the default width
#   * tidy up the sheet view (zoom 110%, select the header row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row content, in final column order -------------------------
$headers = @(
    "BuildingId",        # A
    "BuildingCode",       # B
    "Name",                # C
    "Type",                 # D
    "BuildingClass",         # E
    "ValidFrom",              # F
    "ValidUntil",              # G
    "ConstructionYear",         # H
    "ParkingSpaces",             # I
    "EnergyType",                 # J
    "WaterType",                   # K
    "HeatingType1 ",                # L (trailing space preserved from source)
    "HeatingType2 ",                 # M (trailing space preserved from source)
    "AirConditioning",                # N
    "PrimaryUsage",                    # O
    "SecondaryUsage",                   # P
    "PortfolioId",                       # Q
    "LandId",                             # R
    "SiteId",                              # S
    "CertificateId",                        # T
    "AddressId",                             # U
    "Guid"                                     # V
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Bold the header row (new cellXfs entry with applyFont) ------------
$headerRange = $ws.Range("A1:V1")
$headerRange.Font.Bold = $true

# --- 3. Remove the leftover formatted empty cells at A3 / A5 / A20 --------
$ws.Range("A3").Clear()
$ws.Range("A5").Clear()
$ws.Range("A20").Clear()

# --- 4. Column widths: widen columns whose header no longer fits ----------
# (ColumnWidth = desired xml "width" minus the 5/6 char built-in padding)
$ws.Columns.Item(4).ColumnWidth  = 17.830729166666668  # D  Type            -> 18.6640625
$ws.Columns.Item(5).ColumnWidth  = 12.498697916666666  # E  BuildingClass   -> 13.33203125
$ws.Columns.Item(9).ColumnWidth  = 17.330729166666668  # I  ParkingSpaces   -> 18.1640625
$ws.Columns.Item(11).ColumnWidth = 11.830729166666666  # K  WaterType       -> 12.6640625
$ws.Columns.Item(12).ColumnWidth = 14.166666666666666  # L  HeatingType1    -> 15
$ws.Columns.Item(13).ColumnWidth = 11.166666666666666  # M  HeatingType2    -> 12
$ws.Columns.Item(14).ColumnWidth = 10.830729166666666  # N  AirConditioning -> 11.6640625
$ws.Columns.Item(15).ColumnWidth = 10.330729166666666  # O  PrimaryUsage    -> 11.1640625
$ws.Columns.Item(16).ColumnWidth = 11.830729166666666  # P  SecondaryUsage  -> 12.6640625
$ws.Columns.Item(18).ColumnWidth = 11.330729166666666  # R  LandId          -> 12.1640625
$ws.Columns.Item(19).ColumnWidth = 13.666666666666666  # S  SiteId          -> 14.5

# --- 5. Sheet view tweaks --------------------------------------------------
$excel.ActiveWindow.Zoom = 110
$headerRange.Select()
